$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $ws.Range($cell).Value = "'" + $value
    $ws.Range($cell).Style = "Normal"
}

$ws.Range("D2").Value = "70.353.91"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").Value = "3.790.81"
$ws.Range("E3").Value = "  +0.55%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue "D5" "670.62"
$ws.Range("E5").Value = "  +7.26%  "
Set-TextValue "D6" "168.85"
$ws.Range("E6").Value = "  +1.33%  "
$ws.Range("D7").Value = "3.790.78"
$ws.Range("E7").Value = "  +0.63%  "
Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  -0.02%  "
Set-TextValue "D9" "0.527"
$ws.Range("E9").Value = "  +1.10%  "
Set-TextValue "D10" "0.161"
$ws.Range("E10").Value = "  +0.38%  "
Set-TextValue "D11" "0.463"
$ws.Range("E11").Value = "  +0.97%  "
Set-TextValue "D12" "7.07"
$ws.Range("E12").Value = "  +5.43%  "
Set-TextValue "D13" "0.0000244"
$ws.Range("E13").Value = "  -1.33%  "
Set-TextValue "D14" "35.65"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "4.425.20"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "3.785.55"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "70.266.67"
$ws.Range("E17").Value = "  +1.69%  "
Set-TextValue "D18" "17.64"
$ws.Range("E18").Value = "  -0.03%  "
Set-TextValue "D19" "7.18"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("E20").Value = "  +0.51%  "
Set-TextValue "D21" "11.45"
$ws.Range("E21").Value = "  +19.32%  "
Set-TextValue "D22" "474.14"
$ws.Range("E22").Value = "  +1.26%  "
Set-TextValue "D23" "0.712"
$ws.Range("E23").Value = "  +0.68%  "
Set-TextValue "D24" "82.97"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("E25").Value = "  -4.01%  "
Set-TextValue "D26" "12.23"
$ws.Range("E26").Value = "  +1.27%  "
Set-TextValue "D27" "10.34"
$ws.Range("E27").Value = "  +3.02%  "
Set-TextValue "D28" "2.12"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "3.937.07"
$ws.Range("E30").Value = "  +0.38%  "
Set-TextValue "D31" "2.85"
$ws.Range("E31").Value = "  +6.29%  "
$ws.Range("E32").Value = "  +3.05%  "
Set-TextValue "D33" "7.44"
$ws.Range("E33").Value = "  +2.99%  "
Set-TextValue "D34" "29.51"
$ws.Range("E34").Value = "  +2.65%  "
Set-TextValue "D35" "0.179"
$ws.Range("E35").Value = "  +8.96%  "
$ws.Range("E36").Value = "  +0.09%  "
Set-TextValue "D37" "9.10"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").Value = "3.742.80"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("E39").Value = "  +0.50%  "
Set-TextValue "D40" "3.39"
$ws.Range("E40").Value = "  -1.54%  "
Set-TextValue "D41" "5.96"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("E44").Value = "  +10.75%  "
Set-TextValue "D46" "45.53"
$ws.Range("E46").Value = "  +5.56%  "
Set-TextValue "D47" "158.90"
$ws.Range("E47").Value = "  +4.18%  "
Set-TextValue "D48" "48.03"
$ws.Range("E48").Value = "  +2.97%  "
$ws.Range("E49").Value = "  +4.92%  "
Set-TextValue "D50" "0.300"
$ws.Range("E50").Value = "  +0.75%  "
Set-TextValue "D51" "8.51"
$ws.Range("E51").Value = "  +1.14%  "

# Swap rows 42 and 43 data (B, C, D, E) per the diff
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D42" "0.962"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D43" "0.999"
$ws.Range("E43").Value = "  -0.06%  "
